# Update "想去人数" (want-to-go count) and "最低票价" (min ticket price) figures
# for the latest data refresh (gh-pages output regenerated at 456a3b4).
# The same underlying event rows are duplicated across the "展览" sheet and
# the combined "全部类型" sheet, so both need the identical update.

$wb = $excel.ActiveWorkbook

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1153
    $ws.Range("F3").Value = 586
    $ws.Range("F6").Value = 143

    $ws.Range("F10").Value = 5222
    $ws.Range("G10").Value = 63

    $ws.Range("F11").Value = 4793

    $ws.Range("F13").Value = 35
}
